$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.081.69'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.890.82'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5177'
$ws.Range("E7").Value = '  +2.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3754'
$ws.Range("E8").Value = '  +2.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07211'
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.12'
$ws.Range("E10").Value = '  +1.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9010'
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07654'
$ws.Range("E12").Value = '  +1.81%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.890.55'
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.25'
$ws.Range("E14").Value = '  -1.10%  '
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008512'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.41'
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9998'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.134.88'
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.055'
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.151.64'
$ws.Range("E22").Value = '  +2.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.58'
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.403'
$ws.Range("E24").Value = '  -0.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.292'
$ws.Range("E25").Value = '  +10.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.25'
$ws.Range("E26").Value = '  -1.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.732'
$ws.Range("E27").Value = '  -3.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.05'
$ws.Range("E28").Value = '  +0.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.30'
$ws.Range("E29").Value = '  +0.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.963'
$ws.Range("E30").Value = '  +5.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.795'
$ws.Range("E31").Value = '  +2.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09189'
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05046'
$ws.Range("E33").Value = '  -1.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.237'
$ws.Range("E34").Value = '  +6.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7757'
$ws.Range("E35").Value = '  +3.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.981'
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("E37").Value = '  +1.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.586'
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01989'
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.073'
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.021'
$ws.Range("E42").Value = '  +5.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.36'
$ws.Range("E43").Value = '  +2.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.634'
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1510'
$ws.Range("E45").Value = '  +2.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4828'
$ws.Range("E46").Value = '  +2.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.17'
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.44'
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.87'
$ws.Range("E51").Value = '  +0.97%  '
